$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing question text (A8) to include the difficulty suffix
$ws.Range("A8").Value = "드라마 OST로 안방극장 탐험(난이도 중)"

# Add new description text in B8, matching the style/format of the other description cells
$ws.Range("B7").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("B8").Value = "들려오는 곡이 나온 '드라마의 이름'을 입력해 주세요!"

# Update active selection to B8
$ws.Range("B8").Select()
